$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.824999999999998
$ws.Range("E6").Value = 12.967
$ws.Range("E7").Value = 12.97
$ws.Range("D8").Value = -8.016
$ws.Range("E8").Value = 13.03
$ws.Range("B12").Value = 5.649
$ws.Range("D12").Value = -8.260999999999999
$ws.Range("D14").Value = -8.263
$ws.Range("E19").Value = 12.452
$ws.Range("E21").Value = 13.252
$ws.Range("D22").Value = -8.191999999999998
$ws.Range("E24").Value = 12.81
